$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B (pushes former firstName/lastName columns to C/D)
$ws.Columns.Item(2).Insert()

# Header
$ws.Range("B1").Value = "email"

# Data rows
$ws.Range("B2").Value = "nguyen@gmail.com"
$ws.Range("B3").Value = "dat@gmail.com"

# Turn the email cells into mailto hyperlinks (also applies the built-in
# Hyperlink style / font automatically)
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:nguyen@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:dat@gmail.com")

# Matches the cursor position left behind in the saved file
$null = $ws.Range("X14").Select()
